# Apply the commit's changes to SwaadSutra_Consolidated_2026-01-19.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": update order in row 16 ---
$orders = $wb.Worksheets.Item("All Orders")

# Status column (H) changes from NEW to CANCELLED
$orders.Range("H16").Value = "CANCELLED"

# Cancel Reason column (M) gets a note
$orders.Range("M16").Value = "test order"

# --- Sheet "Daily Summary": update the aggregated stats for 2026-01-13 (row 4) ---
$summary = $wb.Worksheets.Item("Daily Summary")

# Cancelled count increases from 2 to 3
$summary.Range("D4").Value = 3

# Revenue drops from 260 to 200 (the cancelled order's amount removed)
$summary.Range("E4").Value = 200

# Pending amount drops from 260 to 200 as well
$summary.Range("G4").Value = 200
